$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.841467
$ws.Range("H2").Value = 26.524401
$ws.Range("I2").Value = 0.5917001192060068
$ws.Range("J2").Value = 0.5917001192060067
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2226943333333333
$ws.Range("N2").Value = 0.668083
$ws.Range("O2").Value = 0.0947101322715019
$ws.Range("P2").Value = 0.09471013227150192
$ws.Range("Q2").Value = 1.968944599253666
$ws.Range("R2").Value = 17.720501393283
$ws.Range("S2").Value = 0.05603999655506436
$ws.Range("T2").Value = 0.05603999655506435
$ws.Range("G3").Value = 8.841467
$ws.Range("H3").Value = 26.524401
$ws.Range("I3").Value = 0.5917001192060068
$ws.Range("J3").Value = 0.5917001192060067
$ws.Range("O3").Value = 0.3538103900551972
$ws.Range("P3").Value = 0.3538103900551972
$ws.Range("Q3").Value = 7.355422698196666
$ws.Range("R3").Value = 66.19880428377
$ws.Range("S3").Value = 0.2093496499719839
$ws.Range("T3").Value = 0.2093496499719839
$ws.Range("G4").Value = 8.841467
$ws.Range("H4").Value = 26.524401
$ws.Range("I4").Value = 0.5917001192060068
$ws.Range("J4").Value = 0.5917001192060067
$ws.Range("M4").Value = 1.296707666666667
$ws.Range("N4").Value = 3.890123
$ws.Range("O4").Value = 0.5514794776733007
$ws.Range("P4").Value = 0.5514794776733009
$ws.Range("Q4").Value = 11.46479804348033
$ws.Range("R4").Value = 103.183182391323
$ws.Range("S4").Value = 0.3263104726789585
$ws.Range("T4").Value = 0.3263104726789585
$ws.Range("I5").Value = 0.2746155987184545
$ws.Range("J5").Value = 0.2746155987184545
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2226943333333333
$ws.Range("N5").Value = 0.668083
$ws.Range("O5").Value = 0.0947101322715019
$ws.Range("P5").Value = 0.09471013227150192
$ws.Range("Q5").Value = 0.9138123897846666
$ws.Range("R5").Value = 8.224311508061998
$ws.Range("S5").Value = 0.02600887967844252
$ws.Range("T5").Value = 0.02600887967844252
$ws.Range("I6").Value = 0.2746155987184545
$ws.Range("J6").Value = 0.2746155987184545
$ws.Range("O6").Value = 0.3538103900551972
$ws.Range("P6").Value = 0.3538103900551972
$ws.Range("S6").Value = 0.0971618520978179
$ws.Range("T6").Value = 0.0971618520978179
$ws.Range("I7").Value = 0.2746155987184545
$ws.Range("J7").Value = 0.2746155987184545
$ws.Range("M7").Value = 1.296707666666667
$ws.Range("N7").Value = 3.890123
$ws.Range("O7").Value = 0.5514794776733007
$ws.Range("P7").Value = 0.5514794776733009
$ws.Range("Q7").Value = 5.320959514291332
$ws.Range("R7").Value = 47.88863562862199
$ws.Range("S7").Value = 0.1514448669421941
$ws.Range("T7").Value = 0.1514448669421941
$ws.Range("G8").Value = 1.997574666666667
$ws.Range("H8").Value = 5.992724
$ws.Range("I8").Value = 0.1336842820755386
$ws.Range("J8").Value = 0.1336842820755386
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2226943333333333
$ws.Range("N8").Value = 0.668083
$ws.Range("O8").Value = 0.0947101322715019
$ws.Range("P8").Value = 0.09471013227150192
$ws.Range("Q8").Value = 0.4448485586768889
$ws.Range("R8").Value = 4.003637028092
$ws.Range("S8").Value = 0.01266125603799503
$ws.Range("T8").Value = 0.01266125603799503
$ws.Range("G9").Value = 1.997574666666667
$ws.Range("H9").Value = 5.992724
$ws.Range("I9").Value = 0.1336842820755386
$ws.Range("J9").Value = 0.1336842820755386
$ws.Range("O9").Value = 0.3538103900551972
$ws.Range("P9").Value = 0.3538103900551972
$ws.Range("Q9").Value = 1.661828975275556
$ws.Range("R9").Value = 14.95646077748
$ws.Range("S9").Value = 0.04729888798539532
$ws.Range("T9").Value = 0.04729888798539532
$ws.Range("G10").Value = 1.997574666666667
$ws.Range("H10").Value = 5.992724
$ws.Range("I10").Value = 0.1336842820755386
$ws.Range("J10").Value = 0.1336842820755386
$ws.Range("M10").Value = 1.296707666666667
$ws.Range("N10").Value = 3.890123
$ws.Range("O10").Value = 0.5514794776733007
$ws.Range("P10").Value = 0.5514794776733009
$ws.Range("Q10").Value = 2.590270385005778
$ws.Range("R10").Value = 23.312433465052
$ws.Range("S10").Value = 0.07372413805214824
$ws.Range("T10").Value = 0.07372413805214824
